$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "24.402.67"
Set-TextValue $ws.Range("E2") "  -1.43%  "
Set-TextValue $ws.Range("D3") "1.683.37"
Set-TextValue $ws.Range("E3") "  -0.92%  "
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  -0.25%  "
Set-TextValue $ws.Range("D5") "316.11"
Set-TextValue $ws.Range("E5") "  -0.01%  "
Set-TextValue $ws.Range("D6") "0.9996"
Set-TextValue $ws.Range("E6") "  -0.35%  "
Set-TextValue $ws.Range("D7") "0.3884"
Set-TextValue $ws.Range("E7") "  -1.12%  "
Set-TextValue $ws.Range("D8") "0.4006"
Set-TextValue $ws.Range("E8") "  -0.76%  "
Set-TextValue $ws.Range("D9") "1.482"
Set-TextValue $ws.Range("E9") "  -1.40%  "
Set-TextValue $ws.Range("D10") "0.9995"
Set-TextValue $ws.Range("E10") "  -0.38%  "
Set-TextValue $ws.Range("D11") "52.29"
Set-TextValue $ws.Range("E11") "  -3.19%  "
Set-TextValue $ws.Range("D12") "0.08742"
Set-TextValue $ws.Range("E12") "  -1.57%  "
Set-TextValue $ws.Range("D13") "26.12"
Set-TextValue $ws.Range("E13") "  +11.85%  "
Set-TextValue $ws.Range("D14") "7.474"
Set-TextValue $ws.Range("E14") "  +3.46%  "
Set-TextValue $ws.Range("D15") "7.989"
Set-TextValue $ws.Range("E15") "  -0.63%  "
Set-TextValue $ws.Range("D16") "0.00001341"
Set-TextValue $ws.Range("E16") "  +1.07%  "
Set-TextValue $ws.Range("D17") "1.667.72"
Set-TextValue $ws.Range("E17") "  -2.59%  "
Set-TextValue $ws.Range("D18") "97.67"
Set-TextValue $ws.Range("E18") "  -2.41%  "
Set-TextValue $ws.Range("D19") "0.07208"
Set-TextValue $ws.Range("E19") "  +2.88%  "
Set-TextValue $ws.Range("D20") "19.68"
Set-TextValue $ws.Range("E20") "  +0.43%  "
Set-TextValue $ws.Range("D21") "7.250"
Set-TextValue $ws.Range("E21") "  +3.44%  "
Set-TextValue $ws.Range("D22") "1.000"
Set-TextValue $ws.Range("D23") "14.14"
Set-TextValue $ws.Range("E23") "  -2.21%  "
Set-TextValue $ws.Range("D24") "24.406.71"
Set-TextValue $ws.Range("E24") "  -1.35%  "
Set-TextValue $ws.Range("D25") "3.015"
Set-TextValue $ws.Range("E25") "  -7.12%  "
Set-TextValue $ws.Range("D26") "2.338"
Set-TextValue $ws.Range("E26") "  -0.68%  "
Set-TextValue $ws.Range("D27") "22.49"
Set-TextValue $ws.Range("E27") "  -1.02%  "
Set-TextValue $ws.Range("D28") "167.84"
Set-TextValue $ws.Range("E28") "  +4.27%  "
Set-TextValue $ws.Range("D29") "8.584"
Set-TextValue $ws.Range("E29") "  +11.22%  "
Set-TextValue $ws.Range("D30") "5.360"
Set-TextValue $ws.Range("E30") "  +3.89%  "
Set-TextValue $ws.Range("D31") "138.13"
Set-TextValue $ws.Range("E31") "  +1.20%  "
Set-TextValue $ws.Range("D32") "1.853.32"
Set-TextValue $ws.Range("E32") "  -2.33%  "
Set-TextValue $ws.Range("D33") "0.08746"
Set-TextValue $ws.Range("E33") "  +0.09%  "
Set-TextValue $ws.Range("D34") "7.326"
Set-TextValue $ws.Range("E34") "  +2.16%  "
Set-TextValue $ws.Range("D35") "1.045"
Set-TextValue $ws.Range("E35") "  -2.59%  "
Set-TextValue $ws.Range("D36") "0.02997"
Set-TextValue $ws.Range("E36") "  +9.72%  "
Set-TextValue $ws.Range("E37") "  +0.29%  "
Set-TextValue $ws.Range("D38") "0.2752"
Set-TextValue $ws.Range("E38") "  +0.42%  "
Set-TextValue $ws.Range("D39") "10.78"
Set-TextValue $ws.Range("E39") "  -4.12%  "
Set-TextValue $ws.Range("D40") "0.09138"
Set-TextValue $ws.Range("E40") "  -0.35%  "
Set-TextValue $ws.Range("D41") "0.7974"
Set-TextValue $ws.Range("E41") "  +4.13%  "
Set-TextValue $ws.Range("D42") "14.03"
Set-TextValue $ws.Range("E42") "  -2.31%  "
Set-TextValue $ws.Range("D43") "1.469"
Set-TextValue $ws.Range("E43") "  +0.39%  "
Set-TextValue $ws.Range("D44") "17.41"
Set-TextValue $ws.Range("E44") "  +9.62%  "
Set-TextValue $ws.Range("D45") "0.7200"
Set-TextValue $ws.Range("E45") "  +0.55%  "
Set-TextValue $ws.Range("D46") "2.596"
Set-TextValue $ws.Range("E46") "  +0.93%  "
Set-TextValue $ws.Range("D47") "4.260"
Set-TextValue $ws.Range("E47") "  +1.05%  "
Set-TextValue $ws.Range("D48") "1.397"
Set-TextValue $ws.Range("E48") "  +6.58%  "
Set-TextValue $ws.Range("E49") "  -0.34%  "
Set-TextValue $ws.Range("D50") "139.11"
Set-TextValue $ws.Range("E50") "  -1.16%  "
Set-TextValue $ws.Range("D51") "0.08044"
Set-TextValue $ws.Range("E51") "  +0.86%  "
